# Add Singapore surprise songs (rows 156-167): dress colour, song, and mashup columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E156").Value = 'yellow'
$ws.Range("G156").Value = 'Mine (Taylor''s Version)'
$ws.Range("H156").Value = 'Starlight (Taylor''s Version)'

$ws.Range("E157").Value = 'yellow'
$ws.Range("G157").Value = 'I Don''t Wanna Live Forever'
$ws.Range("H157").Value = 'Dress'

$ws.Range("E158").Value = 'green'
$ws.Range("G158").Value = 'long story short'
$ws.Range("H158").Value = 'The Story Of Us (Taylor''s Version)'

$ws.Range("E159").Value = 'green'
$ws.Range("G159").Value = 'Clean (Taylor''s Version)'
$ws.Range("H159").Value = 'evermore'

$ws.Range("E160").Value = 'red'
$ws.Range("G160").Value = 'Foolish One (Taylor''s Version) [From The Vault]'
$ws.Range("H160").Value = 'Tell Me Why (Taylor''s Version)'

$ws.Range("E161").Value = 'red'
$ws.Range("G161").Value = 'This Love (Taylor''s Version)'
$ws.Range("H161").Value = 'Call It What You Want'

$ws.Range("E162").Value = 'red'
$ws.Range("G162").Value = 'Death By A Thousand Cuts'
$ws.Range("H162").Value = 'Babe (Taylor''s Version) [From The Vault]'

$ws.Range("E163").Value = 'red'
$ws.Range("G163").Value = 'Fifteen (Taylor''s Version)'
$ws.Range("H163").Value = 'You''re On Your Own, Kid'

$ws.Range("E164").Value = 'yellow'
$ws.Range("G164").Value = 'Sparks Fly (Taylor''s Version)'
$ws.Range("H164").Value = 'gold rush'

$ws.Range("E165").Value = 'yellow'
$ws.Range("G165").Value = 'False God'
$ws.Range("H165").Value = '"Slut!" (Taylor''s Version) [From The Vault]'

$ws.Range("E166").Value = 'blue'
$ws.Range("G166").Value = 'Tim McGraw'
$ws.Range("H166").Value = 'cowboy like me'

$ws.Range("E167").Value = 'blue'
$ws.Range("G167").Value = 'mirrorball'
$ws.Range("H167").Value = 'epiphany'

# Restore the cursor/selection and scroll position to where the author left it
# (frozen header row stays ySplit=1; the visible window had scrolled to row 141).
$win = $excel.ActiveWindow
$ws.Range("G164").Select()
$win.ScrollRow = 141
$win.ScrollColumn = 1
